# Hawaii Samples Run 11/13/2019
# Add a new data row (65) to Sheet1, mirroring the existing CRM-accuracy rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A64's formatting (date number format) down into A65, then set the date.
$ws.Range("A64").Copy($ws.Range("A65"))
$ws.Range("A65").Value2 = 43782               # 11/13/2019

$ws.Range("B65").Value = 2211.2751980152798
$ws.Range("C65").Value = 2207.0300000000002
$ws.Range("D65").Formula = "=100*(B65-C65)/C65"
$ws.Range("E65").Value = 169
$ws.Range("F65").Value = "crm opened 11/8/2019. capped - no evap"

# Match the shifted active-cell selection left behind by data entry.
$ws.Range("F66").Select()
